$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. Title: "Lucas Hoff Schmidt" -> "CV - Lucas Hoff Schmidt"
Replace-Text "Lucas Hoff Schmidt" "CV - Lucas Hoff Schmidt"

# 2. Profile: "hvordan virksomheder skaber værdi" -> "forretningsudvikling"
Replace-Text "hvordan virksomheder skaber værdi" "forretningsudvikling"

# 3. Profile paragraph rewrite
Replace-Text "I mig får i en kollega der er nytænkende og god til at finde praktiske løsninger, uden at gå på kompromis med kvaliteten. Jeg er en teamspiller med et stærkt selvstændigt drive og " "Jeg er nytænkende og hurtig til at finde praktiske løsninger, uden at gå på kompromis med kvaliteten. Jeg er en holdspiller med et stærkt selvstændigt drive og "

# 4. Section heading: "Evner" -> "Kompetencer"
Replace-Text "Evner" "Kompetencer"

# 5. Skills list: insert "Datarensning • " before SQL, and replace "Machine Learning • Hypothesis Testing •" with
#    "Deskriptiv statistik • Hypotesetestning • Machine Learning •"
Replace-Text "SQL • PostgreSQL" "Datarensning • SQL • PostgreSQL"
Replace-Text "Machine Learning • Hypothesis Testing •" "Deskriptiv statistik • Hypotesetestning • Machine Learning •"

# 6. Systems/platforms bullet list: "Dynamics 365 CRM • REST API • Azure" -> "Dynamics 365 CRM • Azure • Github • REST API"
Replace-Text "Dynamics 365 CRM • REST API • Azure" "Dynamics 365 CRM • Azure • Github • REST API"

# 7. Business/process skills: insert "• Koordinering " after "Procesoptimering "
Replace-Text "Procesoptimering • Formidling" "Procesoptimering • Koordinering • Formidling"

# 8. Education dates: "2017-2023" -> "CBS"
Replace-Text "2017-2023" "CBS"

# 9. Add trailing period to several bullet sentences
Replace-Text "udgav et spil på Google Play med løbende funktionsopdateringer" "udgav et spil på Google Play med løbende funktionsopdateringer."
Replace-Text "Analyserede data for at forbedre spillet og min markedsføringsstrategi" "Analyserede data for at forbedre spillet og min markedsføringsstrategi."
Replace-Text "varemærkeerhvervelse" "varemærkeerhvervelse."

# 10. Insert "diverse " before "udholdenhedsudfordringer."
Replace-Text "angribende tyre og udholdenhedsudfordringer." "angribende tyre og diverse udholdenhedsudfordringer."

# 11. "Rejsen udviklede min problemløsningskompetence," -> "Rejsen udviklede mine problemløsningskompetencer,"
Replace-Text "Rejsen udviklede min problemløsningskompetence," "Rejsen udviklede mine problemløsningskompetencer,"

# 12. More trailing periods
Replace-Text "Validerede fakturaer fra kreditorer med korrektion af uoverensstemmelser" "Validerede fakturaer fra kreditorer med korrektion af uoverensstemmelser."
Replace-Text "Skabte brugermanualer og onboardede nye medarbejdere" "Skabte brugermanualer og onboardede nye medarbejdere."
Replace-Text "Analyserede juridiske dokumenter for at sikre compliance med eksterne love og interne virksomhedspolitikker" "Analyserede juridiske dokumenter for at sikre compliance med eksterne love og interne virksomhedspolitikker."

# 13. "læringspensummer" -> "læringspensa" and add trailing period
Replace-Text "Designede individuelle læringspensummer for at adressere mentees" "Designede individuelle læringspensa for at adressere mentees"
Replace-Text "mentees’ udfordringer" "mentees’ udfordringer."

# 14. Final trailing period
Replace-Text "akademisk og motivationsmæssigt for at hjælpe dem med at opnå deres fulde potentiale" "akademisk og motivationsmæssigt for at hjælpe dem med at opnå deres fulde potentiale."

Write-Output "done"
